# Update the workbook's embedded "version" build string throughout the
# workbook, reflecting the release of "Coal Mine Boundaries and Methane
# Sources - version 1.0.0 (built on February 03 2026 17.29.55 EST)".

$wb = $excel.ActiveWorkbook

$oldVersion = "mines - version 1.0.0 (Feb 3 2026) (built on February 03 2026 10.14.00 EST)"
$newVersion = "Coal Mine Boundaries and Methane Sources - version 1.0.0 (built on February 03 2026 17.29.55 EST)"

# --- "About" sheet ---
$wsAbout = $wb.Worksheets.Item("About")

$wsAbout.Range("A2").Value = "Version: " + $newVersion

$wsAbout.Range("A6").Value = "Recommended Citation:  `"Global Energy Monitor, Coal mine boundaries and methane sources for Tongde Coal Mine, China, M2174, version '" + $newVersion + "'. (See the CC license for attribution requirements if sharing or adapting the data set.)"

# --- "Boundaries and methane sources" sheet ---
$wsData = $wb.Worksheets.Item("Boundaries and methane sources")

foreach ($r in 2..7) {
    $cell = $wsData.Cells.Item($r, 19)  # column S = 19
    if ($cell.Text -eq $oldVersion) {
        $cell.Value = $newVersion
    }
}
